$wb = $excel.ActiveWorkbook

# --- Sheet ALC ---
$ws = $wb.Worksheets.Item("ALC")
# row 98
$ws.Range("H98").Value = 902.56525
$ws.Range("I98").Value = 927.1667
$ws.Range("J98").Value = 814
$ws.Range("K98").Value = 927.1667
$ws.Range("L98").Value = 814
$ws.Range("M98").Value = 570.8333
$ws.Range("N98").Value = -3810
# row 122
$ws.Range("H122").Value = 902.56525
$ws.Range("I122").Value = 927.1667
$ws.Range("J122").Value = 814
$ws.Range("K122").Value = 2781.5001
$ws.Range("L122").Value = 2442
$ws.Range("M122").Value = -331.5001000000002
$ws.Range("N122").Value = -7342
# row 124
$ws.Range("H124").Value = 78150
$ws.Range("J124").Value = 78150
$ws.Range("L124").Value = 78150
$ws.Range("N124").Value = -87970
# row 138
$ws.Range("H138").Value = 2785.822
$ws.Range("I138").Value = 870.7143
$ws.Range("J138").Value = 4549.737
$ws.Range("K138").Value = 2612.1429
$ws.Range("L138").Value = 13649.211
$ws.Range("M138").Value = 2527.8571
$ws.Range("N138").Value = -23929.211

# --- Sheet ARM ---
$ws = $wb.Worksheets.Item("ARM")
# row 8
$ws.Range("H8").Value = 342600
$ws.Range("I8").Value = 504000
$ws.Range("J8").Value = 19800
$ws.Range("K8").Value = 504000
$ws.Range("L8").Value = 19800
$ws.Range("M8").Value = -503856
$ws.Range("N8").Value = -20088
# row 11
$ws.Range("H11").Value = 1000000
$ws.Range("I11").Value = 1000000
$ws.Range("J11").Value = 0
$ws.Range("K11").Value = 1000000
$ws.Range("L11").Value = 0
$ws.Range("M11").Value = -999856
$ws.Range("N11").ClearContents()
# row 12
$ws.Range("H12").Value = 0
$ws.Range("I12").Value = 0
$ws.Range("K12").Value = 0
$ws.Range("M12").ClearContents()
# row 13
$ws.Range("H13").Value = 255975
$ws.Range("J13").Value = 7966.6665
$ws.Range("L13").Value = 7966.6665
$ws.Range("N13").Value = -8254.666499999999
# row 14
$ws.Range("H14").Value = 14800
$ws.Range("I14").Value = 0
$ws.Range("J14").Value = 14800
$ws.Range("K14").Value = 0
$ws.Range("L14").Value = 14800
$ws.Range("M14").ClearContents()
$ws.Range("N14").Value = -15150
# row 32
$ws.Range("H32").Value = 17753.854
$ws.Range("I32").Value = 19724.807
$ws.Range("J32").Value = 8353.923000000001
$ws.Range("K32").Value = 19724.807
$ws.Range("L32").Value = 8353.923000000001
$ws.Range("M32").Value = -19437.807
$ws.Range("N32").Value = -8927.923000000001
# row 61
$ws.Range("H61").Value = 8970.5
$ws.Range("I61").Value = 6358.1953
$ws.Range("J61").Value = 17209.309
$ws.Range("K61").Value = 6358.1953
$ws.Range("L61").Value = 17209.309
$ws.Range("M61").Value = -6146.1953
$ws.Range("N61").Value = -17633.309
# row 122
$ws.Range("H122").Value = 2976855.2
$ws.Range("I122").Value = 622.95
$ws.Range("J122").Value = 62501500
$ws.Range("K122").Value = 1868.85
$ws.Range("L122").Value = 187504500
$ws.Range("M122").Value = 581.1499999999999
$ws.Range("N122").Value = -187509400
# row 132
$ws.Range("H132").Value = 3687.2322
$ws.Range("I132").Value = 1195.4186
$ws.Range("K132").Value = 3586.2558
$ws.Range("M132").Value = -1056.2558
# row 136
$ws.Range("H136").Value = 8970.5
$ws.Range("I136").Value = 6358.1953
$ws.Range("J136").Value = 17209.309
$ws.Range("K136").Value = 19074.5859
$ws.Range("L136").Value = 51627.927
$ws.Range("M136").Value = -16524.5859
$ws.Range("N136").Value = -56727.927

# --- Sheet BSM ---
$ws = $wb.Worksheets.Item("BSM")
# row 16
$ws.Range("H16").Value = 18790
$ws.Range("J16").Value = 18790
$ws.Range("L16").Value = 18790
$ws.Range("N16").Value = -19130
# row 103
$ws.Range("H103").Value = 40000
$ws.Range("J103").Value = 40000
$ws.Range("L103").Value = 40000
$ws.Range("N103").Value = -42344

# --- Sheet CRP ---
$ws = $wb.Worksheets.Item("CRP")
# row 31
$ws.Range("H31").Value = 3668.2666
$ws.Range("I31").Value = 0
$ws.Range("J31").Value = 3668.2666
$ws.Range("K31").Value = 0
$ws.Range("L31").Value = 3668.2666
$ws.Range("M31").ClearContents()
$ws.Range("N31").Value = -4258.2666
# row 34
$ws.Range("H34").Value = 3668.2666
$ws.Range("I34").Value = 0
$ws.Range("J34").Value = 3668.2666
$ws.Range("K34").Value = 0
$ws.Range("L34").Value = 3668.2666
$ws.Range("M34").ClearContents()
$ws.Range("N34").Value = -4072.2666
# row 132
$ws.Range("H132").Value = 5177.275
$ws.Range("I132").Value = 6493.864
$ws.Range("J132").Value = 3568.111
$ws.Range("K132").Value = 19481.592
$ws.Range("L132").Value = 10704.333
$ws.Range("M132").Value = -16951.592
$ws.Range("N132").Value = -15764.333
# row 134
$ws.Range("H134").Value = 3148.1042
$ws.Range("I134").Value = 2362.158
$ws.Range("K134").Value = 7086.474
$ws.Range("M134").Value = -4551.474
# row 140
$ws.Range("H140").Value = 54300
$ws.Range("J140").Value = 54300
$ws.Range("L140").Value = 54300
$ws.Range("N140").Value = -64660

# --- Sheet CUL ---
$ws = $wb.Worksheets.Item("CUL")
# row 2
$ws.Range("H2").Value = 44.121212
$ws.Range("I2").Value = 10
$ws.Range("J2").Value = 48.827587
$ws.Range("K2").Value = 60
$ws.Range("L2").Value = 292.965522
$ws.Range("M2").Value = 53
$ws.Range("N2").Value = -518.965522
# row 117
$ws.Range("H117").Value = 37037180
$ws.Range("I117").Value = 214.5
$ws.Range("J117").Value = 111111110
$ws.Range("K117").Value = 643.5
$ws.Range("L117").Value = 333333330
$ws.Range("M117").Value = 2798.5
$ws.Range("N117").Value = -333340214
# row 132
$ws.Range("H132").Value = 2499.9167
$ws.Range("I132").Value = 4416.6665
$ws.Range("J132").Value = 1861
$ws.Range("K132").Value = 39749.9985
$ws.Range("L132").Value = 16749
$ws.Range("M132").Value = -37219.9985
$ws.Range("N132").Value = -21809

# --- Sheet GSM ---
$ws = $wb.Worksheets.Item("GSM")
# row 132
$ws.Range("H132").Value = 3472.2964
$ws.Range("I132").Value = 1623.6586
$ws.Range("J132").Value = 9302.615
$ws.Range("K132").Value = 4870.9758
$ws.Range("L132").Value = 27907.845
$ws.Range("M132").Value = -2340.9758
$ws.Range("N132").Value = -32967.845

# --- Sheet LTW ---
$ws = $wb.Worksheets.Item("LTW")
# row 61
$ws.Range("H61").Value = 724424.4399999999
$ws.Range("I61").Value = 16994.572
$ws.Range("J61").Value = 1431854.2
$ws.Range("K61").Value = 16994.572
$ws.Range("L61").Value = 1431854.2
$ws.Range("M61").Value = -16792.572
$ws.Range("N61").Value = -1432258.2
# row 113
$ws.Range("H113").Value = 724424.4399999999
$ws.Range("I113").Value = 16994.572
$ws.Range("J113").Value = 1431854.2
$ws.Range("K113").Value = 16994.572
$ws.Range("L113").Value = 1431854.2
$ws.Range("M113").Value = -14824.572
$ws.Range("N113").Value = -1436194.2
# row 132
$ws.Range("H132").Value = 4298.8037
$ws.Range("I132").Value = 4242.615
$ws.Range("J132").Value = 4481.4165
$ws.Range("K132").Value = 12727.845
$ws.Range("L132").Value = 13444.2495
$ws.Range("M132").Value = -10197.845
$ws.Range("N132").Value = -18504.2495
# row 136
$ws.Range("H136").Value = 4685.75
$ws.Range("I136").Value = 2754.7188
$ws.Range("K136").Value = 8264.1564
$ws.Range("M136").Value = -5714.1564

# --- Sheet WVR ---
$ws = $wb.Worksheets.Item("WVR")
# row 49
$ws.Range("H49").Value = 30262
$ws.Range("J49").Value = 30262
$ws.Range("L49").Value = 30262
$ws.Range("N49").Value = -30722
# row 54
$ws.Range("H54").Value = 15912.786
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 15912.786
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 15912.786
$ws.Range("M54").ClearContents()
$ws.Range("N54").Value = -16952.786
# row 122
$ws.Range("H122").Value = 2148.36
$ws.Range("I122").Value = 2091.318
$ws.Range("J122").Value = 2566.6667
$ws.Range("K122").Value = 6273.954000000001
$ws.Range("L122").Value = 7700.000100000001
$ws.Range("M122").Value = -3823.954000000001
$ws.Range("N122").Value = -12600.0001
# row 125
$ws.Range("H125").Value = 61143.332
$ws.Range("J125").Value = 61143.332
$ws.Range("L125").Value = 61143.332
$ws.Range("N125").Value = -70983.33199999999
# row 132
$ws.Range("H132").Value = 1801.7805
$ws.Range("I132").Value = 888.56525
$ws.Range("K132").Value = 2665.69575
$ws.Range("M132").Value = -135.6957499999999
